# Add a new "Year" column (column E) to the surface meteorological stations sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column D (style index "1") onto the new column E
# so the new cells pick up the same left-aligned style used by the rest
# of the table.
$ws.Range("D1:D108").Copy()
$ws.Range("E1:E108").PasteSpecial(-4122)

# Header for the new column.
$ws.Range("E1").Value = "Year"

# Every station record was collected in 2022.
$ws.Range("E2:E108").Value = 2022

# Match the selection left behind in the authored workbook.
$ws.Range("E2:E108").Select()
